$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.277799999999997
$ws.Range("D14").Value = -8.146400000000002
$ws.Range("D21").Value = -7.766099999999999
$ws.Range("D23").Value = -6.762399999999993
$ws.Range("D25").Value = -8.394499999999999
